# ShiftGuard_System.xlsx - "Fix formula references and add comprehensive documentation"
#
# 1) Settings sheet: the two placeholder/spacer rows (28 and 33) still carried
#    two empty inline-string cells each (leftover from the original template).
#    Clear them out so the rows collapse down to bare, contentless rows -
#    matching the cleanup of the other blank spacer rows in that sheet.
# 2) Calculations!J2 and Payroll!G2/H2/I2/J2/L2/M2/N2 were off-by-one (and
#    more) against the Settings sheet layout - they pointed at section
#    *header* rows (or otherwise wrong rows) instead of the actual values
#    beneath them. Re-point each at the correct Settings!$B$xx cell.

$wb = $excel.ActiveWorkbook

# --- 1) Settings: tidy up the now-pointless empty cells in the blank spacer rows ---
$settings = $wb.Worksheets.Item("Settings")

$settings.Rows.Item(28).OutlineLevel = 0
$settings.Cells.Item(28, 1).ClearContents()
$settings.Cells.Item(28, 2).ClearContents()

$settings.Rows.Item(33).OutlineLevel = 0
$settings.Cells.Item(33, 1).ClearContents()
$settings.Cells.Item(33, 2).ClearContents()

# --- 2) Calculations sheet: fix Attendance Bonus reference ---
$calculations = $wb.Worksheets.Item("Calculations")
$calculations.Range("J2").Formula = '=IF(AND(BiometricData!A2<>"",G2>=100),Settings!$B$23,"")'

# --- 3) Payroll sheet: fix Base Pay / Overtime Pay / Night Allowance / Weekend
#        Allowance / Tax Deduction / Health Insurance / Pension Deduction refs ---
$payroll = $wb.Worksheets.Item("Payroll")
$payroll.Range("G2").Formula = '=IF(A2<>"",D2*Settings!$B$23,"")'
$payroll.Range("H2").Formula = '=IF(A2<>"",E2*Settings!$B$23*Settings!$B$25,"")'
$payroll.Range("I2").Formula = '=IF(A2<>"",F2*Settings!$B$23*Settings!$B$26/100,"")'
$payroll.Range("J2").Formula = '=IF(AND(A2<>"",OR(WEEKDAY(B2)=1,WEEKDAY(B2)=7)),D2*Settings!$B$23*Settings!$B$27/100,0)'
$payroll.Range("L2").Formula = '=IF(A2<>"",K2*Settings!$B$30/100,"")'
$payroll.Range("M2").Formula = '=IF(A2<>"",Settings!$B$31,"")'
$payroll.Range("N2").Formula = '=IF(A2<>"",K2*Settings!$B$32/100,"")'
